$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 53 (shifts existing rows 53:121 down to 54:122)
$ws.Rows(53).Insert()

# The table ("Cluster_Keywords") originally spanned A1:C121; after the sheet-level
# row insert it still reports A1:C121 (the newly-inserted row sits just above the
# old last row, which slid outside the table range). Resize it back to include
# the now-empty row 53 and the shifted last row (now at 122).
$tbl = $ws.ListObjects.Item("Cluster_Keywords")
$tbl.Resize($ws.Range("A1:C122"))

# Fill in the values for the newly inserted row.
$ws.Range("A53").Value = "Poult"
$ws.Range("B53").Formula = "=LEN(Cluster_Keywords[[#This Row],[Stem]])"
$ws.Range("C53").Value = "Food & Drink"

# Make sure the calculated column formula on the row that got shifted to the
# bottom of the table (row 122) keeps using the table structured reference
# form, rather than any simplified alias.
$ws.Range("B122").Formula = "=LEN(Cluster_Keywords[[#This Row],[Stem]])"

# Reflect the cell that was active/selected in the sheet when the edit was saved.
[void]$ws.Range("C54").Select()
